# Update Kujata market price/profit data (H:N columns) across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4633879
$ws.Range("I62").Value = 6177137.5
$ws.Range("J62").Value = 4102.6665
$ws.Range("K62").Value = 6177137.5
$ws.Range("L62").Value = 4102.6665
$ws.Range("M62").Value = -6176513.5
$ws.Range("N62").Value = -5350.6665
$ws.Range("H65").Value = 4633879
$ws.Range("I65").Value = 6177137.5
$ws.Range("J65").Value = 4102.6665
$ws.Range("K65").Value = 30885687.5
$ws.Range("L65").Value = 20513.3325
$ws.Range("M65").Value = -30882567.5
$ws.Range("N65").Value = -26753.3325
$ws.Range("H106").Value = 6077.7144
$ws.Range("I106").Value = 6280.7407
$ws.Range("K106").Value = 6280.7407
$ws.Range("M106").Value = -5649.7407
$ws.Range("H129").Value = 948.1
$ws.Range("J129").Value = 999.1429000000001
$ws.Range("L129").Value = 2997.4287
$ws.Range("N129").Value = -12997.4287
$ws.Range("H132").Value = 6067667
$ws.Range("I132").Value = 6807545
$ws.Range("K132").Value = 20422635
$ws.Range("M132").Value = -20420105
$ws.Range("H138").Value = 3537.3494
$ws.Range("I138").Value = 2017.5
$ws.Range("J138").Value = 3958.2307
$ws.Range("K138").Value = 6052.5
$ws.Range("L138").Value = 11874.6921
$ws.Range("M138").Value = -912.5
$ws.Range("N138").Value = -22154.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2193.5625
$ws.Range("I2").Value = 1465.75
$ws.Range("K2").Value = 1465.75
$ws.Range("M2").Value = -1352.75
$ws.Range("H32").Value = 19428.031
$ws.Range("I32").Value = 17127.518
$ws.Range("K32").Value = 17127.518
$ws.Range("M32").Value = -16840.518
$ws.Range("H61").Value = 40001540
$ws.Range("I61").Value = 47620370
$ws.Range("K61").Value = 47620370
$ws.Range("M61").Value = -47620158
$ws.Range("H74").Value = 1810.174
$ws.Range("I74").Value = 709.5625
$ws.Range("K74").Value = 709.5625
$ws.Range("M74").Value = 164.4375
$ws.Range("H77").Value = 1810.174
$ws.Range("I77").Value = 709.5625
$ws.Range("K77").Value = 3547.8125
$ws.Range("M77").Value = 820.1875
$ws.Range("H116").Value = 2193.5625
$ws.Range("I116").Value = 1465.75
$ws.Range("K116").Value = 1465.75
$ws.Range("M116").Value = 828.25
$ws.Range("H132").Value = 3143.9062
$ws.Range("I132").Value = 2318.2354
$ws.Range("J132").Value = 4079.6667
$ws.Range("K132").Value = 6954.706200000001
$ws.Range("L132").Value = 12239.0001
$ws.Range("M132").Value = -4424.706200000001
$ws.Range("N132").Value = -17299.0001
$ws.Range("H136").Value = 40001540
$ws.Range("I136").Value = 47620370
$ws.Range("K136").Value = 142861110
$ws.Range("M136").Value = -142858560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2193.5625
$ws.Range("I3").Value = 1465.75
$ws.Range("K3").Value = 1465.75
$ws.Range("M3").Value = -1351.75
$ws.Range("H105").Value = 76924800
$ws.Range("I105").Value = 100001544
$ws.Range("J105").Value = 2333.3333
$ws.Range("K105").Value = 100001544
$ws.Range("L105").Value = 2333.3333
$ws.Range("M105").Value = -99999797
$ws.Range("N105").Value = -5827.3333
$ws.Range("H107").Value = 1037.55
$ws.Range("I107").Value = 976.6667
$ws.Range("J107").Value = 1220.2
$ws.Range("K107").Value = 976.6667
$ws.Range("L107").Value = 1220.2
$ws.Range("M107").Value = 943.3333
$ws.Range("N107").Value = -5060.2
$ws.Range("H134").Value = 3620.738
$ws.Range("I134").Value = 959.9167
$ws.Range("J134").Value = 7168.5
$ws.Range("K134").Value = 2879.7501
$ws.Range("L134").Value = 21505.5
$ws.Range("M134").Value = -344.7501000000002
$ws.Range("N134").Value = -26575.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1582.5902
$ws.Range("I31").Value = 1391.9286
$ws.Range("J31").Value = 3718
$ws.Range("K31").Value = 1391.9286
$ws.Range("L31").Value = 3718
$ws.Range("M31").Value = -1096.9286
$ws.Range("N31").Value = -4308
$ws.Range("H34").Value = 1582.5902
$ws.Range("I34").Value = 1391.9286
$ws.Range("J34").Value = 3718
$ws.Range("K34").Value = 1391.9286
$ws.Range("M34").Value = -1189.9286
$ws.Range("N34").Value = -4122
$ws.Range("H86").Value = 2104545.5
$ws.Range("I86").Value = 2792612.5
$ws.Range("J86").Value = 40344.875
$ws.Range("K86").Value = 2792612.5
$ws.Range("L86").Value = 40344.875
$ws.Range("M86").Value = -2791489.5
$ws.Range("N86").Value = -42590.875
$ws.Range("H89").Value = 2104545.5
$ws.Range("I89").Value = 2792612.5
$ws.Range("J89").Value = 40344.875
$ws.Range("K89").Value = 13963062.5
$ws.Range("L89").Value = 201724.375
$ws.Range("M89").Value = -13957446.5
$ws.Range("N89").Value = -212956.375
$ws.Range("H105").Value = 863.4545000000001
$ws.Range("I105").Value = 833.1111
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 833.1111
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 913.8889
$ws.Range("N105").Value = -4494
$ws.Range("H132").Value = 2379.258
$ws.Range("I132").Value = 2097.3044
$ws.Range("J132").Value = 3189.875
$ws.Range("K132").Value = 6291.9132
$ws.Range("L132").Value = 9569.625
$ws.Range("M132").Value = -3761.9132
$ws.Range("N132").Value = -14629.625
$ws.Range("H134").Value = 12501650
$ws.Range("I134").Value = 1586.0646
$ws.Range("J134").Value = 55557424
$ws.Range("K134").Value = 4758.1938
$ws.Range("L134").Value = 166672272
$ws.Range("M134").Value = -2223.1938
$ws.Range("N134").Value = -166677342
$ws.Range("H141").Value = 410099.12
$ws.Range("J141").Value = 410099.12
$ws.Range("L141").Value = 410099.12
$ws.Range("N141").Value = -420459.12

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 283.55554
$ws.Range("I50").Value = 168.33333
$ws.Range("J50").Value = 341.16666
$ws.Range("K50").Value = 504.99999
$ws.Range("L50").Value = 1023.49998
$ws.Range("M50").Value = -23.99998999999997
$ws.Range("N50").Value = -1985.49998
$ws.Range("H53").Value = 283.55554
$ws.Range("I53").Value = 168.33333
$ws.Range("J53").Value = 341.16666
$ws.Range("K53").Value = 504.99999
$ws.Range("L53").Value = 1023.49998
$ws.Range("M53").Value = -23.99998999999997
$ws.Range("N53").Value = -1985.49998
$ws.Range("H131").Value = 24427898
$ws.Range("J131").Value = 41667.137
$ws.Range("L131").Value = 125001.411
$ws.Range("N131").Value = -135081.411

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 6000
$ws.Range("J20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("N20").Value = -6490
$ws.Range("H24").Value = 5000
$ws.Range("J24").Value = 5000
$ws.Range("L24").Value = 5000
$ws.Range("N24").Value = -5346
$ws.Range("H97").Value = 621.8570999999999
$ws.Range("J97").Value = 645.44446
$ws.Range("L97").Value = 645.44446
$ws.Range("N97").Value = -1637.44446
$ws.Range("H132").Value = 4735.5
$ws.Range("I132").Value = 5190.6
$ws.Range("K132").Value = 15571.8
$ws.Range("M132").Value = -13041.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1993.4286
$ws.Range("I136").Value = 1778.5625
$ws.Range("J136").Value = 2681
$ws.Range("K136").Value = 5335.6875
$ws.Range("L136").Value = 8043
$ws.Range("M136").Value = -2785.6875
$ws.Range("N136").Value = -13143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 421.72726
$ws.Range("I100").Value = 409.875
$ws.Range("J100").Value = 453.33334
$ws.Range("K100").Value = 819.75
$ws.Range("L100").Value = 906.66668
$ws.Range("M100").Value = -278.75
$ws.Range("N100").Value = -1988.66668
$ws.Range("H132").Value = 3832.9268
$ws.Range("I132").Value = 3969.3103
$ws.Range("J132").Value = 3503.3333
$ws.Range("K132").Value = 11907.9309
$ws.Range("L132").Value = 10509.9999
$ws.Range("M132").Value = -9377.930899999999
$ws.Range("N132").Value = -15569.9999
$ws.Range("H136").Value = 1402.8235
$ws.Range("I136").Value = 494.8
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 1484.4
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = 1065.6
$ws.Range("N136").Value = -13200
